$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.937.30"
$ws.Range("D3").Value = "1.817.16"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.42"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07369"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8706"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "1.812.79"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.383"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07115"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.512"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.36"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008682"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.65"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "26.967.80"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("D24").Value = "2.056.77"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.894"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.13"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.258"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.29"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08879"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7595"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.896"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.094"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05287"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01946"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.980"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5289"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.327"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.429"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4853"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.31"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.660"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06291"
$ws.Range("E51").Value = "  -0.04%  "